$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Servicios" label moves from J1 to M1
$ws.Range("J1").ClearContents()
$ws.Range("M1").Value = "Servicios"

# New column header "Empresa Prestadora" in H2
$ws.Range("H2").Value = "Empresa Prestadora"

# Remove the old "Servicio" label (column J) entirely - no replacement
$ws.Range("J2").ClearContents()

# Old Servicios list (J3:J5) moves to M3:M5, with a new numbering column L3:L6
# and a new 4th service "Testeo CoVid" in M6
$ws.Range("J3").ClearContents()
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = "Atencion Medica"

$ws.Range("J4").ClearContents()
$ws.Range("L4").Value = 2
$ws.Range("M4").Value = "Traslado"

$ws.Range("J5").ClearContents()
$ws.Range("L5").Value = 3
$ws.Range("M5").Value = "Urgencias"

# New note about "Empresa Prestadora" filtering
$ws.Range("A24").Value = "La lista de servicios es por empresa prestadora. Cuando un propsecto va a pasar a cliente le filtro la combo por empresa prestadora cargada."

$ws.Range("L6").Value = 4
$ws.Range("M6").Value = "Testeo CoVid"

# New "ServiciosPorPlan" style mapping rows for "Azul PAMI"
$ws.Range("O3").Value = "Azul PAMI"
$ws.Range("P3").Value = 1
$ws.Range("O4").Value = "Azul PAMI"
$ws.Range("P4").Value = 4

# Updated price values for "Azul Urgencias VIP" (row 4) and "Azul Clientes" (row 5)
$ws.Range("B4").Value = 100
$ws.Range("C4").Value = 150
$ws.Range("D4").Value = 300
$ws.Range("E4").Value = 400

$ws.Range("B5").Value = 220
$ws.Range("C5").Value = 280
$ws.Range("D5").Value = 400
$ws.Range("E5").Value = 600

# Column width adjustments for the newly used columns
$ws.Columns.Item(8).ColumnWidth = 20.1666666667
$ws.Columns.Item(13).ColumnWidth = 16.33
$ws.Columns.Item(15).ColumnWidth = 16.0

# Restore the active cell selection as last left by the author
$ws.Range("Q32").Select()
